# Daily attendance processing - 2026-01-24 05:37:50
# For every cell in column G ("Recorded By") that contains a two-part,
# comma-separated list of recorders (e.g. "System, someone@example.com"),
# reverse the order of the two parts (e.g. "someone@example.com, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value2

    if ($null -ne $val -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val.Split(",")
        if ($parts.Count -eq 2) {
            $first = $parts[0].Trim()
            $second = $parts[1].Trim()
            $cell.Value2 = "$second, $first"
        }
    }
}
